# e-Prescribe app.pptx edit script
# - Fills in the (previously blank) "Resources" slide (5th slide) with a
#   title and a credit line for the icons used.
# - Inserts a brand-new blank slide (Title+Content layout) right after the
#   Resources slide - this becomes the new slide with id 263, pushing the
#   "Struggles" / "Learnings" slides one position later.

$p = $ppt.ActivePresentation

# --- 1) Duplicate the still-blank "Resources" slide (slide 5) first, so the
#        duplicate ends up blank too, and lands immediately after it (slide 6).
$resourcesSlide = $p.Slides.Item(5)
$resourcesSlide.Duplicate() | Out-Null

# --- 2) Now fill in the original slide 5 ("Resources") with its title and
#        the icons credit line in the content placeholder.
$titleRange = $resourcesSlide.Shapes.Item(1).TextFrame.TextRange
$titleRun = $titleRange.InsertAfter("Resources")
$titleRun.LanguageID = "en-PH"

$bodyShape = $resourcesSlide.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRun = $bodyRange.InsertAfter("Icons - https://icons8.com/")
$bodyRun.LanguageID = "en-PH"

# Remove the bullet and flush the paragraph to the left margin (no indent),
# matching the plain "Icons - https://icons8.com/" credit-line formatting.
$bodyRange.ParagraphFormat.Bullet.Type = 0
$bodyRuler = $bodyShape.TextFrame.Ruler.Levels.Item(1)
$bodyRuler.FirstMargin = 0
